$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: apply corrected match ordering for several rounds ---
# Source data for betexplorer had home/away rows in the wrong order for a
# handful of fixture pairs/groups that share a matchday; this restores the
# correct per-row association between rows and betting data (content for
# columns F:V moves between rows while Indice/pais/torneio/temporada/data_partida
# in A:E stay put).

# Row 2 <- content sourced from original row 4
$ws.Cells.Item(2, 6).Value = "Vrsac"
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = "Kolubara"
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 2.71
$ws.Cells.Item(2, 11).Value = "05/08/2023 13:12"
$ws.Cells.Item(2, 12).Value = 2.71
$ws.Cells.Item(2, 13).Value = "05/08/2023 13:12"
$ws.Cells.Item(2, 14).Value = 2.81
$ws.Cells.Item(2, 15).Value = "05/08/2023 13:12"
$ws.Cells.Item(2, 16).Value = 2.83
$ws.Cells.Item(2, 17).Value = "05/08/2023 15:33"
$ws.Cells.Item(2, 18).Value = 2.53
$ws.Cells.Item(2, 19).Value = "05/08/2023 13:12"
$ws.Cells.Item(2, 20).Value = 2.53
$ws.Cells.Item(2, 21).Value = "05/08/2023 13:12"
$ws.Cells.Item(2, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/vrsac-kolubara/rVQRcbII/"

# Row 4 <- content sourced from original row 2
$ws.Cells.Item(4, 6).Value = "Smederevo"
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = "FK Indjija"
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 2.74
$ws.Cells.Item(4, 11).Value = "05/08/2023 13:12"
$ws.Cells.Item(4, 12).Value = 2.74
$ws.Cells.Item(4, 13).Value = "05/08/2023 13:12"
$ws.Cells.Item(4, 14).Value = 2.75
$ws.Cells.Item(4, 15).Value = "05/08/2023 13:12"
$ws.Cells.Item(4, 16).Value = 2.76
$ws.Cells.Item(4, 17).Value = "05/08/2023 15:33"
$ws.Cells.Item(4, 18).Value = 2.56
$ws.Cells.Item(4, 19).Value = "05/08/2023 13:12"
$ws.Cells.Item(4, 20).Value = 2.56
$ws.Cells.Item(4, 21).Value = "05/08/2023 13:12"
$ws.Cells.Item(4, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/smederevo-indjija/rR4gggd8/"

# Row 5 <- content sourced from original row 6
$ws.Cells.Item(5, 6).Value = "Sloboda"
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = "OFK Beograd"
$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 10).Value = 1.79
$ws.Cells.Item(5, 11).Value = "05/08/2023 13:12"
$ws.Cells.Item(5, 12).Value = 1.89
$ws.Cells.Item(5, 13).Value = "05/08/2023 15:40"
$ws.Cells.Item(5, 14).Value = 3.1
$ws.Cells.Item(5, 15).Value = "05/08/2023 13:12"
$ws.Cells.Item(5, 16).Value = 3.03
$ws.Cells.Item(5, 17).Value = "05/08/2023 18:02"
$ws.Cells.Item(5, 18).Value = 4.17
$ws.Cells.Item(5, 19).Value = "05/08/2023 13:12"
$ws.Cells.Item(5, 20).Value = 3.96
$ws.Cells.Item(5, 21).Value = "05/08/2023 15:40"
$ws.Cells.Item(5, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/sloboda-ofk-beograd/Iqc6jitR/"

# Row 6 <- content sourced from original row 5
$ws.Cells.Item(6, 6).Value = "Metalac"
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = "Graficar Beograd"
$ws.Cells.Item(6, 9).Value = 3
$ws.Cells.Item(6, 10).Value = 2.01
$ws.Cells.Item(6, 11).Value = "04/08/2023 08:12"
$ws.Cells.Item(6, 12).Value = 2.17
$ws.Cells.Item(6, 13).Value = "04/08/2023 12:33"
$ws.Cells.Item(6, 14).Value = 3
$ws.Cells.Item(6, 15).Value = "04/08/2023 08:12"
$ws.Cells.Item(6, 16).Value = 3.11
$ws.Cells.Item(6, 17).Value = "05/08/2023 18:02"
$ws.Cells.Item(6, 18).Value = 3.15
$ws.Cells.Item(6, 19).Value = "04/08/2023 08:12"
$ws.Cells.Item(6, 20).Value = 3.06
$ws.Cells.Item(6, 21).Value = "04/08/2023 12:33"
$ws.Cells.Item(6, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/metalac-graficar-beograd/zPqYoBlr/"

# Row 25 <- content sourced from original row 26
$ws.Cells.Item(25, 6).Value = "FK Indjija"
$ws.Cells.Item(25, 7).Value = 2
$ws.Cells.Item(25, 8).Value = "RFK Novi Sad"
$ws.Cells.Item(25, 9).Value = 1
$ws.Cells.Item(25, 10).Value = 1.45
$ws.Cells.Item(25, 11).Value = "25/08/2023 08:13"
$ws.Cells.Item(25, 12).Value = 1.34
$ws.Cells.Item(25, 13).Value = "26/08/2023 16:45"
$ws.Cells.Item(25, 14).Value = 3.74
$ws.Cells.Item(25, 15).Value = "25/08/2023 08:13"
$ws.Cells.Item(25, 16).Value = 4.28
$ws.Cells.Item(25, 17).Value = "26/08/2023 16:45"
$ws.Cells.Item(25, 18).Value = 5.1
$ws.Cells.Item(25, 19).Value = "25/08/2023 08:13"
$ws.Cells.Item(25, 20).Value = 8.12
$ws.Cells.Item(25, 21).Value = "26/08/2023 16:45"
$ws.Cells.Item(25, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/indjija-rfk-novi-sad/OKWOf6XC/"

# Row 26 <- content sourced from original row 25
$ws.Cells.Item(26, 6).Value = "Smederevo"
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = "Macva"
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 1.81
$ws.Cells.Item(26, 11).Value = "26/08/2023 13:43"
$ws.Cells.Item(26, 12).Value = 2.25
$ws.Cells.Item(26, 13).Value = "26/08/2023 16:55"
$ws.Cells.Item(26, 14).Value = 3.15
$ws.Cells.Item(26, 15).Value = "26/08/2023 13:43"
$ws.Cells.Item(26, 16).Value = 2.92
$ws.Cells.Item(26, 17).Value = "26/08/2023 16:55"
$ws.Cells.Item(26, 18).Value = 3.97
$ws.Cells.Item(26, 19).Value = "26/08/2023 13:43"
$ws.Cells.Item(26, 20).Value = 3.14
$ws.Cells.Item(26, 21).Value = "26/08/2023 16:55"
$ws.Cells.Item(26, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/smederevo-macva-sabac/l6zmaLOb/"

# Row 30 <- content sourced from original row 31
$ws.Cells.Item(30, 6).Value = "Graficar Beograd"
$ws.Cells.Item(30, 7).Value = 1
$ws.Cells.Item(30, 8).Value = "Dubocica"
$ws.Cells.Item(30, 9).Value = 2
$ws.Cells.Item(30, 10).Value = 1.69
$ws.Cells.Item(30, 11).Value = "27/08/2023 11:43"
$ws.Cells.Item(30, 12).Value = 1.73
$ws.Cells.Item(30, 13).Value = "27/08/2023 16:49"
$ws.Cells.Item(30, 14).Value = 3.3
$ws.Cells.Item(30, 15).Value = "27/08/2023 11:43"
$ws.Cells.Item(30, 16).Value = 3.29
$ws.Cells.Item(30, 17).Value = "27/08/2023 16:49"
$ws.Cells.Item(30, 18).Value = 4.38
$ws.Cells.Item(30, 19).Value = "27/08/2023 11:43"
$ws.Cells.Item(30, 20).Value = 4.45
$ws.Cells.Item(30, 21).Value = "27/08/2023 16:49"
$ws.Cells.Item(30, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/graficar-beograd-dubocica/2BVSgQnJ/"

# Row 31 <- content sourced from original row 30
$ws.Cells.Item(31, 6).Value = "OFK Beograd"
$ws.Cells.Item(31, 7).Value = 2
$ws.Cells.Item(31, 8).Value = "Radnicki Beograd"
$ws.Cells.Item(31, 9).Value = 2
$ws.Cells.Item(31, 10).Value = 1.68
$ws.Cells.Item(31, 11).Value = "27/08/2023 11:43"
$ws.Cells.Item(31, 12).Value = 1.49
$ws.Cells.Item(31, 13).Value = "27/08/2023 16:58"
$ws.Cells.Item(31, 14).Value = 3.29
$ws.Cells.Item(31, 15).Value = "27/08/2023 11:43"
$ws.Cells.Item(31, 16).Value = 3.84
$ws.Cells.Item(31, 17).Value = "27/08/2023 16:58"
$ws.Cells.Item(31, 18).Value = 4.48
$ws.Cells.Item(31, 19).Value = "27/08/2023 11:43"
$ws.Cells.Item(31, 20).Value = 5.78
$ws.Cells.Item(31, 21).Value = "27/08/2023 16:58"
$ws.Cells.Item(31, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-radnicki-beograd/Emrz2Nvt/"

# Row 42 <- content sourced from original row 44
$ws.Cells.Item(42, 6).Value = "FK Indjija"
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = "Graficar Beograd"
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 1.96
$ws.Cells.Item(42, 11).Value = "08/09/2023 04:42"
$ws.Cells.Item(42, 12).Value = 2.09
$ws.Cells.Item(42, 13).Value = "09/09/2023 13:14"
$ws.Cells.Item(42, 14).Value = 3.12
$ws.Cells.Item(42, 15).Value = "08/09/2023 04:42"
$ws.Cells.Item(42, 16).Value = 3.02
$ws.Cells.Item(42, 17).Value = "09/09/2023 14:31"
$ws.Cells.Item(42, 18).Value = 3.16
$ws.Cells.Item(42, 19).Value = "08/09/2023 04:42"
$ws.Cells.Item(42, 20).Value = 3.33
$ws.Cells.Item(42, 21).Value = "09/09/2023 13:14"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/indjija-graficar-beograd/xjeQCvgb/"

# Row 43 <- content sourced from original row 42
$ws.Cells.Item(43, 6).Value = "Smederevo"
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(43, 8).Value = "Metalac"
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 10).Value = 2.3
$ws.Cells.Item(43, 11).Value = "09/09/2023 13:42"
$ws.Cells.Item(43, 12).Value = 2.44
$ws.Cells.Item(43, 13).Value = "09/09/2023 16:21"
$ws.Cells.Item(43, 14).Value = 2.76
$ws.Cells.Item(43, 15).Value = "09/09/2023 13:42"
$ws.Cells.Item(43, 16).Value = 2.69
$ws.Cells.Item(43, 17).Value = "09/09/2023 15:59"
$ws.Cells.Item(43, 18).Value = 3.18
$ws.Cells.Item(43, 19).Value = "09/09/2023 13:42"
$ws.Cells.Item(43, 20).Value = 3.1
$ws.Cells.Item(43, 21).Value = "09/09/2023 16:21"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/smederevo-metalac/vFis9dhN/"

# Row 44 <- content sourced from original row 43
$ws.Cells.Item(44, 6).Value = "Radnicki S. Mitrovica"
$ws.Cells.Item(44, 7).Value = 4
$ws.Cells.Item(44, 8).Value = "Jedinstvo U."
$ws.Cells.Item(44, 9).Value = 1
$ws.Cells.Item(44, 10).Value = 2.04
$ws.Cells.Item(44, 11).Value = "08/09/2023 04:42"
$ws.Cells.Item(44, 12).Value = 2.33
$ws.Cells.Item(44, 13).Value = "09/09/2023 16:27"
$ws.Cells.Item(44, 14).Value = 2.88
$ws.Cells.Item(44, 15).Value = "08/09/2023 04:42"
$ws.Cells.Item(44, 16).Value = 3.03
$ws.Cells.Item(44, 17).Value = "09/09/2023 15:46"
$ws.Cells.Item(44, 18).Value = 3.23
$ws.Cells.Item(44, 19).Value = "08/09/2023 04:42"
$ws.Cells.Item(44, 20).Value = 2.89
$ws.Cells.Item(44, 21).Value = "09/09/2023 16:27"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-s-mitrovica-jedinstvo-ub/AumwAxwH/"

# Row 49 <- content sourced from original row 50
$ws.Cells.Item(49, 6).Value = "Sloboda"
$ws.Cells.Item(49, 7).Value = 1
$ws.Cells.Item(49, 8).Value = "Vrsac"
$ws.Cells.Item(49, 9).Value = 1
$ws.Cells.Item(49, 10).Value = 2.2
$ws.Cells.Item(49, 11).Value = "15/09/2023 06:12"
$ws.Cells.Item(49, 12).Value = 2.5
$ws.Cells.Item(49, 13).Value = "16/09/2023 18:52"
$ws.Cells.Item(49, 14).Value = 2.71
$ws.Cells.Item(49, 15).Value = "15/09/2023 06:12"
$ws.Cells.Item(49, 16).Value = 2.62
$ws.Cells.Item(49, 17).Value = "16/09/2023 18:45"
$ws.Cells.Item(49, 18).Value = 3.18
$ws.Cells.Item(49, 19).Value = "15/09/2023 06:12"
$ws.Cells.Item(49, 20).Value = 3.09
$ws.Cells.Item(49, 21).Value = "16/09/2023 18:52"
$ws.Cells.Item(49, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/sloboda-vrsac/jgLaRw8i/"

# Row 50 <- content sourced from original row 49
$ws.Cells.Item(50, 6).Value = "Macva"
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(50, 8).Value = "FK Indjija"
$ws.Cells.Item(50, 9).Value = 1
$ws.Cells.Item(50, 10).Value = 2.33
$ws.Cells.Item(50, 11).Value = "15/09/2023 06:12"
$ws.Cells.Item(50, 12).Value = 2.9
$ws.Cells.Item(50, 13).Value = "16/09/2023 18:57"
$ws.Cells.Item(50, 14).Value = 2.71
$ws.Cells.Item(50, 15).Value = "15/09/2023 06:12"
$ws.Cells.Item(50, 16).Value = 2.88
$ws.Cells.Item(50, 17).Value = "16/09/2023 18:57"
$ws.Cells.Item(50, 18).Value = 2.87
$ws.Cells.Item(50, 19).Value = "15/09/2023 06:12"
$ws.Cells.Item(50, 20).Value = 2.42
$ws.Cells.Item(50, 21).Value = "16/09/2023 18:57"
$ws.Cells.Item(50, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/macva-sabac-indjija/xjXH2fVj/"

# Row 52 <- content sourced from original row 53
$ws.Cells.Item(52, 6).Value = "RFK Novi Sad"
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = "Tekstilac Odzaci"
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 2.88
$ws.Cells.Item(52, 11).Value = "16/09/2023 03:13"
$ws.Cells.Item(52, 12).Value = 4.03
$ws.Cells.Item(52, 13).Value = "17/09/2023 15:39"
$ws.Cells.Item(52, 14).Value = 2.86
$ws.Cells.Item(52, 15).Value = "16/09/2023 03:13"
$ws.Cells.Item(52, 16).Value = 3.11
$ws.Cells.Item(52, 17).Value = "17/09/2023 15:03"
$ws.Cells.Item(52, 18).Value = 2.22
$ws.Cells.Item(52, 19).Value = "16/09/2023 03:13"
$ws.Cells.Item(52, 20).Value = 1.86
$ws.Cells.Item(52, 21).Value = "17/09/2023 15:39"
$ws.Cells.Item(52, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/rfk-novi-sad-tekstilac-odzaci/UJl8Ne7G/"

# Row 53 <- content sourced from original row 54
$ws.Cells.Item(53, 6).Value = "Radnicki Beograd"
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = "Radnicki S. Mitrovica"
$ws.Cells.Item(53, 9).Value = 1
$ws.Cells.Item(53, 10).Value = 2.87
$ws.Cells.Item(53, 11).Value = "16/09/2023 03:13"
$ws.Cells.Item(53, 12).Value = 2.49
$ws.Cells.Item(53, 13).Value = "17/09/2023 15:50"
$ws.Cells.Item(53, 14).Value = 2.82
$ws.Cells.Item(53, 15).Value = "16/09/2023 03:13"
$ws.Cells.Item(53, 16).Value = 2.73
$ws.Cells.Item(53, 17).Value = "17/09/2023 15:50"
$ws.Cells.Item(53, 18).Value = 2.25
$ws.Cells.Item(53, 19).Value = "16/09/2023 03:13"
$ws.Cells.Item(53, 20).Value = 2.96
$ws.Cells.Item(53, 21).Value = "17/09/2023 15:50"
$ws.Cells.Item(53, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-beograd-radnicki-s-mitrovica/zck4OyhA/"

# Row 54 <- content sourced from original row 52
$ws.Cells.Item(54, 6).Value = "Graficar Beograd"
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = "OFK Beograd"
$ws.Cells.Item(54, 9).Value = 5
$ws.Cells.Item(54, 10).Value = 2.33
$ws.Cells.Item(54, 11).Value = "16/09/2023 03:13"
$ws.Cells.Item(54, 12).Value = 2.56
$ws.Cells.Item(54, 13).Value = "17/09/2023 15:55"
$ws.Cells.Item(54, 14).Value = 3.02
$ws.Cells.Item(54, 15).Value = "16/09/2023 03:13"
$ws.Cells.Item(54, 16).Value = 3.5
$ws.Cells.Item(54, 17).Value = "17/09/2023 15:55"
$ws.Cells.Item(54, 18).Value = 2.59
$ws.Cells.Item(54, 19).Value = "16/09/2023 03:13"
$ws.Cells.Item(54, 20).Value = 2.11
$ws.Cells.Item(54, 21).Value = "17/09/2023 15:55"
$ws.Cells.Item(54, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/graficar-beograd-ofk-beograd/rNhCMFMM/"

# Row 79 <- content sourced from original row 80
$ws.Cells.Item(79, 6).Value = "Sloboda"
$ws.Cells.Item(79, 7).Value = 3
$ws.Cells.Item(79, 8).Value = "Radnicki Beograd"
$ws.Cells.Item(79, 9).Value = 1
$ws.Cells.Item(79, 10).Value = 1.69
$ws.Cells.Item(79, 11).Value = "13/10/2023 02:13"
$ws.Cells.Item(79, 12).Value = 1.83
$ws.Cells.Item(79, 13).Value = "14/10/2023 14:09"
$ws.Cells.Item(79, 14).Value = 3.03
$ws.Cells.Item(79, 15).Value = "13/10/2023 02:13"
$ws.Cells.Item(79, 16).Value = 3.01
$ws.Cells.Item(79, 17).Value = "14/10/2023 14:09"
$ws.Cells.Item(79, 18).Value = 4.34
$ws.Cells.Item(79, 19).Value = "13/10/2023 02:13"
$ws.Cells.Item(79, 20).Value = 4.44
$ws.Cells.Item(79, 21).Value = "14/10/2023 13:59"
$ws.Cells.Item(79, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/sloboda-radnicki-beograd/YDOkOmFO/"

# Row 80 <- content sourced from original row 79
$ws.Cells.Item(80, 6).Value = "FK Indjija"
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = "Kolubara"
$ws.Cells.Item(80, 9).Value = 2
$ws.Cells.Item(80, 10).Value = 2.03
$ws.Cells.Item(80, 11).Value = "13/10/2023 02:13"
$ws.Cells.Item(80, 12).Value = 2.15
$ws.Cells.Item(80, 13).Value = "14/10/2023 14:51"
$ws.Cells.Item(80, 14).Value = 2.83
$ws.Cells.Item(80, 15).Value = "13/10/2023 02:13"
$ws.Cells.Item(80, 16).Value = 2.81
$ws.Cells.Item(80, 17).Value = "14/10/2023 14:51"
$ws.Cells.Item(80, 18).Value = 3.3
$ws.Cells.Item(80, 19).Value = "13/10/2023 02:13"
$ws.Cells.Item(80, 20).Value = 3.53
$ws.Cells.Item(80, 21).Value = "14/10/2023 14:51"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/indjija-kolubara/CAvE9eHU/"

# Row 97 <- content sourced from original row 100
$ws.Cells.Item(97, 6).Value = "Jedinstvo U."
$ws.Cells.Item(97, 7).Value = 2
$ws.Cells.Item(97, 8).Value = "Radnicki Beograd"
$ws.Cells.Item(97, 9).Value = 1
$ws.Cells.Item(97, 10).Value = 1.52
$ws.Cells.Item(97, 11).Value = "27/10/2023 02:13"
$ws.Cells.Item(97, 12).Value = 1.52
$ws.Cells.Item(97, 13).Value = "28/10/2023 13:57"
$ws.Cells.Item(97, 14).Value = 3.45
$ws.Cells.Item(97, 15).Value = "27/10/2023 02:13"
$ws.Cells.Item(97, 16).Value = 3.63
$ws.Cells.Item(97, 17).Value = "28/10/2023 13:57"
$ws.Cells.Item(97, 18).Value = 4.9
$ws.Cells.Item(97, 19).Value = "27/10/2023 02:13"
$ws.Cells.Item(97, 20).Value = 5.76
$ws.Cells.Item(97, 21).Value = "28/10/2023 13:57"
$ws.Cells.Item(97, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/jedinstvo-ub-radnicki-beograd/lpBfhST4/"

# Row 98 <- content sourced from original row 97
$ws.Cells.Item(98, 6).Value = "Smederevo"
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = "Kolubara"
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 2.29
$ws.Cells.Item(98, 11).Value = "27/10/2023 02:13"
$ws.Cells.Item(98, 12).Value = 2.73
$ws.Cells.Item(98, 13).Value = "28/10/2023 13:58"
$ws.Cells.Item(98, 14).Value = 2.72
$ws.Cells.Item(98, 15).Value = "27/10/2023 02:13"
$ws.Cells.Item(98, 16).Value = 2.6
$ws.Cells.Item(98, 17).Value = "28/10/2023 13:38"
$ws.Cells.Item(98, 18).Value = 2.92
$ws.Cells.Item(98, 19).Value = "27/10/2023 02:13"
$ws.Cells.Item(98, 20).Value = 2.83
$ws.Cells.Item(98, 21).Value = "28/10/2023 13:58"
$ws.Cells.Item(98, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/smederevo-kolubara/nHYMTb2n/"

# Row 100 <- content sourced from original row 101
$ws.Cells.Item(100, 6).Value = "Mladost GAT"
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = "Tekstilac Odzaci"
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 2.09
$ws.Cells.Item(100, 11).Value = "27/10/2023 02:13"
$ws.Cells.Item(100, 12).Value = 2.18
$ws.Cells.Item(100, 13).Value = "28/10/2023 13:51"
$ws.Cells.Item(100, 14).Value = 2.75
$ws.Cells.Item(100, 15).Value = "27/10/2023 02:13"
$ws.Cells.Item(100, 16).Value = 2.9
$ws.Cells.Item(100, 17).Value = "28/10/2023 13:51"
$ws.Cells.Item(100, 18).Value = 3.27
$ws.Cells.Item(100, 19).Value = "27/10/2023 02:13"
$ws.Cells.Item(100, 20).Value = 3.33
$ws.Cells.Item(100, 21).Value = "28/10/2023 13:51"
$ws.Cells.Item(100, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/mladost-gat-tekstilac-odzaci/KUnXzRbU/"

# Row 101 <- content sourced from original row 98
$ws.Cells.Item(101, 6).Value = "Sloboda"
$ws.Cells.Item(101, 7).Value = 2
$ws.Cells.Item(101, 8).Value = "RFK Novi Sad"
$ws.Cells.Item(101, 9).Value = 1
$ws.Cells.Item(101, 10).Value = 1.6
$ws.Cells.Item(101, 11).Value = "27/10/2023 02:13"
$ws.Cells.Item(101, 12).Value = 1.58
$ws.Cells.Item(101, 13).Value = "28/10/2023 13:50"
$ws.Cells.Item(101, 14).Value = 3.24
$ws.Cells.Item(101, 15).Value = "27/10/2023 02:13"
$ws.Cells.Item(101, 16).Value = 3.4
$ws.Cells.Item(101, 17).Value = "28/10/2023 13:50"
$ws.Cells.Item(101, 18).Value = 4.67
$ws.Cells.Item(101, 19).Value = "27/10/2023 02:13"
$ws.Cells.Item(101, 20).Value = 5.63
$ws.Cells.Item(101, 21).Value = "28/10/2023 13:50"
$ws.Cells.Item(101, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/sloboda-rfk-novi-sad/hIHoflah/"

# --- Part 2: append 3 newly scraped fixtures (rows 108-110) ---
# Copy formats (styles only) from the last existing data row (107) down
# onto the new rows, then populate values.
$ws.Range("A107:V107").Copy()
$ws.Range("A108:V110").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 108
$ws.Cells.Item(108, 1).Value = 107
$ws.Cells.Item(108, 2).Value = "serbia"
$ws.Cells.Item(108, 3).Value = "prva-liga"
$ws.Cells.Item(108, 4).Value = "2023-2024"
$ws.Cells.Item(108, 5).Value = 45236.54166666666
$ws.Cells.Item(108, 6).Value = "Radnicki S. Mitrovica"
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = "Smederevo"
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 1.92
$ws.Cells.Item(108, 11).Value = "06/11/2023 01:12"
$ws.Cells.Item(108, 12).Value = 1.78
$ws.Cells.Item(108, 13).Value = "06/11/2023 12:25"
$ws.Cells.Item(108, 14).Value = 2.95
$ws.Cells.Item(108, 15).Value = "06/11/2023 01:12"
$ws.Cells.Item(108, 16).Value = 3.21
$ws.Cells.Item(108, 17).Value = "06/11/2023 12:31"
$ws.Cells.Item(108, 18).Value = 3.83
$ws.Cells.Item(108, 19).Value = "06/11/2023 01:12"
$ws.Cells.Item(108, 20).Value = 4.23
$ws.Cells.Item(108, 21).Value = "06/11/2023 12:25"
$ws.Cells.Item(108, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-s-mitrovica-smederevo/CtiO77Ti/"

# Row 109
$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 2).Value = "serbia"
$ws.Cells.Item(109, 3).Value = "prva-liga"
$ws.Cells.Item(109, 4).Value = "2023-2024"
$ws.Cells.Item(109, 5).Value = 45236.54166666666
$ws.Cells.Item(109, 6).Value = "Radnicki Beograd"
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = "FK Indjija"
$ws.Cells.Item(109, 9).Value = 1
$ws.Cells.Item(109, 10).Value = 2.65
$ws.Cells.Item(109, 11).Value = "06/11/2023 01:12"
$ws.Cells.Item(109, 12).Value = 2.75
$ws.Cells.Item(109, 13).Value = "06/11/2023 12:57"
$ws.Cells.Item(109, 14).Value = 2.86
$ws.Cells.Item(109, 15).Value = "06/11/2023 01:12"
$ws.Cells.Item(109, 16).Value = 2.72
$ws.Cells.Item(109, 17).Value = "06/11/2023 12:58"
$ws.Cells.Item(109, 18).Value = 2.54
$ws.Cells.Item(109, 19).Value = "06/11/2023 01:12"
$ws.Cells.Item(109, 20).Value = 2.68
$ws.Cells.Item(109, 21).Value = "06/11/2023 12:58"
$ws.Cells.Item(109, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-beograd-indjija/SzAbinqB/"

# Row 110
$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 2).Value = "serbia"
$ws.Cells.Item(110, 3).Value = "prva-liga"
$ws.Cells.Item(110, 4).Value = "2023-2024"
$ws.Cells.Item(110, 5).Value = 45236.54166666666
$ws.Cells.Item(110, 6).Value = "RFK Novi Sad"
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = "Metalac"
$ws.Cells.Item(110, 9).Value = 1
$ws.Cells.Item(110, 10).Value = 1.99
$ws.Cells.Item(110, 11).Value = "06/11/2023 01:12"
$ws.Cells.Item(110, 12).Value = 2.31
$ws.Cells.Item(110, 13).Value = "06/11/2023 12:53"
$ws.Cells.Item(110, 14).Value = 2.9
$ws.Cells.Item(110, 15).Value = "06/11/2023 01:12"
$ws.Cells.Item(110, 16).Value = 3.01
$ws.Cells.Item(110, 17).Value = "06/11/2023 12:50"
$ws.Cells.Item(110, 18).Value = 3.76
$ws.Cells.Item(110, 19).Value = "06/11/2023 01:12"
$ws.Cells.Item(110, 20).Value = 2.95
$ws.Cells.Item(110, 21).Value = "06/11/2023 12:53"
$ws.Cells.Item(110, 22).Value = "https://www.betexplorer.com/football/serbia/prva-liga/rfk-novi-sad-metalac/02N5kQDN/"
